$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 268; this shifts rows 268-337 down to 269-338
# (same pattern/format as the existing rows is inherited from the row below).
$ws.Rows.Item(268).Insert()

# Populate the new row 268 with the new record's data.
$ws.Range("A268").Value = 5
$ws.Range("B268").Value = "Macroferia Regional de Talca"
$ws.Range("C268").Value = "Maule"
$ws.Range("D268").Value = 44736
$ws.Range("E268").Value = 7
$ws.Range("F268").Value = 100114014
$ws.Range("G268").Value = "Betarraga"
$ws.Range("H268").Value = "Sin especificar"
$ws.Range("I268").Value = "Primera"
$ws.Range("J268").Value = 5000
$ws.Range("K268").Value = 700
$ws.Range("L268").Value = 700
$ws.Range("M268").Value = 700
$ws.Range("N268").Value = "$/paquete 5 unidades"
$ws.Range("O268").Value = "Región del Maule"
$ws.Range("P268").Value = 140
$ws.Range("Q268").Value = 5
$ws.Range("R268").Value = "Hortaliza"
